# Append: 2025-10-17 12:47 JST
# Update the "取得日時" (acquisition timestamp) column A for the data rows
# on the "ランサーズ" sheet from the previous run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-17 12:47:40"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
